$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '40.345.57'
$ws.Range("E2").Value = '  -3.21%  '
$ws.Range("D3").Value = '2.369.66'
$ws.Range("E3").Value = '  -4.36%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '''310.48'
$ws.Range("E5").Value = '  -2.81%  '
$ws.Range("D6").Value = '''86.19'
$ws.Range("E6").Value = '  -6.99%  '
$ws.Range("D7").Value = '''0.534'
$ws.Range("E7").Value = '  -3.18%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").Value = '''0.490'
$ws.Range("E9").Value = '  -4.24%  '
$ws.Range("D10").Value = '''0.0830'
$ws.Range("E10").Value = '  -4.30%  '
$ws.Range("D11").Value = '''30.38'
$ws.Range("E11").Value = '  -8.42%  '
$ws.Range("D12").Value = '''0.110'
$ws.Range("E12").Value = '  -0.60%  '
$ws.Range("D13").Value = '2.730.45'
$ws.Range("E13").Value = '  -4.57%  '
$ws.Range("D14").Value = '''6.48'
$ws.Range("E14").Value = '  -6.05%  '
$ws.Range("D15").Value = '''14.98'
$ws.Range("E15").Value = '  -3.82%  '
$ws.Range("D16").Value = '2.371.50'
$ws.Range("E16").Value = '  -3.72%  '
$ws.Range("D17").Value = '''0.758'
$ws.Range("E17").Value = '  -4.57%  '
$ws.Range("D18").Value = '40.349.41'
$ws.Range("E18").Value = '  -3.07%  '
$ws.Range("D19").Value = '0.0₃0909'
$ws.Range("E19").Value = '  -3.58%  '
$ws.Range("D20").Value = '''6.13'
$ws.Range("E20").Value = '  -5.13%  '
$ws.Range("D21").Value = '''68.31'
$ws.Range("E21").Value = '  -3.59%  '
$ws.Range("D22").Value = '''10.81'
$ws.Range("E22").Value = '  -4.00%  '
$ws.Range("D23").Value = '''235.13'
$ws.Range("E23").Value = '  -1.95%  '
$ws.Range("D24").Value = '''2.58'
$ws.Range("E24").Value = '  -6.31%  '
$ws.Range("D25").Value = '''1.00'
$ws.Range("E25").Value = '  +0.07%  '
$ws.Range("D26").Value = '''1.81'
$ws.Range("E26").Value = '  -7.32%  '
$ws.Range("D27").Value = '''23.71'
$ws.Range("E27").Value = '  -5.12%  '
$ws.Range("D28").Value = '''2.15'
$ws.Range("E28").Value = '  -3.96%  '
$ws.Range("D29").Value = '''9.24'
$ws.Range("E29").Value = '  -5.23%  '
$ws.Range("D30").Value = '''34.26'
$ws.Range("E30").Value = '  -6.74%  '
$ws.Range("D31").Value = '''152.94'
$ws.Range("E31").Value = '  -3.03%  '
$ws.Range("E32").Value = '  -0.18%  '
$ws.Range("D33").Value = '''5.21'
$ws.Range("E33").Value = '  -4.26%  '
$ws.Range("D34").Value = '''0.0731'
$ws.Range("E34").Value = '  -4.52%  '
$ws.Range("E35").Value = '  -5.54%  '
$ws.Range("E36").Value = '  -2.01%  '
$ws.Range("D37").Value = '''2.79'
$ws.Range("E37").Value = '  -3.32%  '
$ws.Range("D38").Value = '''15.96'
$ws.Range("E38").Value = '  -7.14%  '
$ws.Range("D39").Value = '''0.0998'
$ws.Range("E39").Value = '  -4.08%  '
$ws.Range("D40").Value = '''1.71'
$ws.Range("E40").Value = '  -7.47%  '
$ws.Range("D41").Value = '''3.83'
$ws.Range("E41").Value = '  -4.66%  '
$ws.Range("E42").Value = '  -3.68%  '
$ws.Range("D43").Value = '1.968.73'
$ws.Range("E43").Value = '  -1.49%  '
$ws.Range("D44").Value = '''0.0267'
$ws.Range("E44").Value = '  -5.87%  '
$ws.Range("D45").Value = '''17.74'
$ws.Range("E45").Value = '  -5.06%  '
$ws.Range("D46").Value = '''9.29'
$ws.Range("E46").Value = '  -2.91%  '
$ws.Range("D47").Value = '''2.69'
$ws.Range("E47").Value = '  -9.86%  '
$ws.Range("D48").Value = '2.598.85'
$ws.Range("E48").Value = '  -4.42%  '
$ws.Range("D49").Value = '''93.12'
$ws.Range("E49").Value = '  -5.05%  '
$ws.Range("D50").Value = '''71.42'
$ws.Range("E50").Value = '  -5.58%  '
$ws.Range("D51").Value = '''50.33'
$ws.Range("E51").Value = '  -3.88%  '
